# test_parameters_chl.xlsx cleanup
# Sheets: transitions, costs, utilities, specification, condensed_states
$wb = $excel.ActiveWorkbook

$wsCosts = $wb.Worksheets.Item("costs")
$wsUtil  = $wb.Worksheets.Item("utilities")

# --- costs sheet: insert a "type" column (B) filled with "static" ---
$wsCosts.Columns("B:B").Insert()
$wsCosts.Range("B1").Value = "type"
$wsCosts.Range("B2:B19").Value = "static"

# --- utilities sheet: insert a "type" column (B) filled with "static" ---
$wsUtil.Columns("B:B").Insert()
$wsUtil.Range("B1").Value = "type"
$wsUtil.Range("B2:B19").Value = "static"

# --- restore per-sheet selections ---
$wb.Worksheets.Item("transitions").Range("D5").Select()
$wb.Worksheets.Item("specification").Range("B2").Select()
$wb.Worksheets.Item("condensed_states").Range("G10").Select()
$wsUtil.Range("I14").Select()

# costs ends up the active/selected sheet
$wsCosts.Range("G15").Select()
